$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.253.28'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '3.492.93'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'586.94"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = "'134.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('D7').Value = '3.492.31'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.486"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').Value = '4.084.77'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '3.490.22'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '64.298.10'
$ws.Range('D18').Value = "'25.29"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.25%  '
$ws.Range('E19').Value = '  -1.89%  '
$ws.Range('D20').Value = "'5.75"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = "'13.64"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.80%  '
$ws.Range('D22').Value = "'388.77"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').Value = '3.631.45'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = "'0.565"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.36%  '
$ws.Range('D25').Value = "'74.33"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = "'5.69"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('E31').Value = '  -5.01%  '
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').Value = '3.512.90'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').Value = "'5.23"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').Value = "'6.86"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('D41').Value = "'162.31"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.83%  '
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('D43').Value = "'0.805"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'25.48"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.76%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'1.00"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = "'41.78"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('D50').Value = '2.471.01'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('E51').Value = '  -2.22%  '
